$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N, shifting N:P -> O:Q
$null = $ws.Columns("N").Insert()

# Make "Repayment schedule" the active/selected sheet (was "Prepay Loan")
$ws.Activate()

# Update the selected cell on the newly active sheet
$null = $ws.Range("L9").Select()
